$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Standardize the District column (G3:G30) to the official name.
$ws.Range("G3:G30").Value = "Kalaburagi (Gulbarga)"

# Row 23 (SREEDEVIKA P R) had a stray empty inline-string cell in F23 that
# should be removed entirely (no address data for this entry).
$ws.Range("F23").ClearContents()
